$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-9: Beteckning (A), Datum (B), Forandrad (C), Area (G)
# All "Forandrad" (C) dates advance from 46059 to 46060, and several rows
# are re-ordered to reflect updated record positions.
$data = @(
    @{ Row = 2; A = "A 51680-2024"; B = 45607; C = 46060; G = 1.1 },
    @{ Row = 3; A = "A 34310-2024"; B = 45524; C = 46060; G = 4.8 },
    @{ Row = 4; A = "A 25617-2024"; B = 45463; C = 46060; G = 2.3 },
    @{ Row = 5; A = "A 45983-2023"; B = 45196; C = 46060; G = 0.6 },
    @{ Row = 6; A = "A 54782-2022"; B = 44883; C = 46060; G = 5.5 },
    @{ Row = 7; A = "A 843-2024";   B = 45300; C = 46060; G = 0.8 },
    @{ Row = 8; A = "A 844-2024";   B = 45300; C = 46060; G = 1.2 },
    @{ Row = 9; A = "A 17908-2021"; B = 44301; C = 46060; G = 0.9 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 7).Value = $item.G
}
